$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Row 4: the existing-vehicle filter cell now targets ICE cars instead of BEVs
$ws.Range("AJ4").Value = "T-CAR-ICE*"

# Row 5 (new): Fix O&M cost line for existing ICE cars
$ws.Range("D5").Value = "FIXOM"

$ws.Range("H5").Value = 1
$ws.Range("I5:AH5").Value = 1
$ws.Range("AH5").Style = "Normal"

$ws.Range("AJ5").Value = "T-CAR-ICE*"
$ws.Range("AK5").Value = "*Existing*"

# Restore the view to the state captured in the saved workbook
$ws.Activate() | Out-Null
$ws.Range("AJ14").Select() | Out-Null
